# Config workbook update: add "DIGDAT URL" / "DIGDAT access expiry" rows
# and an "Attachment Path" row to the Assets settings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$lo = $ws.ListObjects.Item("Table1")

# Insert two new blank rows right after "LSBUD URL" (row 11), pushing the
# remaining settings ("Default end date days" downward, etc.) down by two.
[void]$ws.Rows.Item(12).Insert()
[void]$ws.Rows.Item(13).Insert()

# Populate the two new rows (Name column + mirrored placeholder Value column).
$ws.Range("A12").Value = "DIGDAT URL"
$ws.Range("B12").Value = "DIGDAT URL"

$ws.Range("A13").Value = "DIGDAT access expiry"
$ws.Range("B13").Value = "DIGDAT access expiry"

# The old "Search Report Path" row (now shifted down to row 22) is replaced
# with a new "Attachment Path" setting; "Retry Number" stays the last row (23).
$ws.Range("A22").Value = "Attachment Path"
$ws.Range("B22").Value = "Attachment Path"

# Grow the table definition to cover the two newly-inserted rows.
[void]$lo.Resize($ws.Range("A1:B23"))

# Leave the selection where the user last edited.
[void]$ws.Range("B22").Select()
